$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81 (this shifts existing rows 81..158 down to 82..159,
# carrying their formatting/values with them, matching the target diff where every
# row from the old row 81 onward is pushed down by one and a brand-new record is
# inserted at row 81).
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new record's data.
$ws.Range("A81").Value = 4
$ws.Range("B81").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C81").Value = "Los Lagos"
$ws.Range("D81").Value = 44484
$ws.Range("E81").Value = 10
$ws.Range("F81").Value = 100112043
$ws.Range("G81").Value = "Pepino ensalada"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 400
$ws.Range("K81").Value = 17000
$ws.Range("L81").Value = 17000
$ws.Range("M81").Value = 17000
$ws.Range("N81").Value = '$/caja 60 unidades'
$ws.Range("O81").Value = "Región de Arica y Parinacota"
$ws.Range("P81").Value = 283
$ws.Range("Q81").Value = 60
$ws.Range("R81").Value = "Hortaliza"
